$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with the new block order labels
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "living_rooms_2"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "bedrooms_1"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Update the one-hot block-order matrix rows 2-7
$values = @(
  @(0,0,0,0,1,0),
  @(0,0,1,0,0,0),
  @(1,0,0,0,0,0),
  @(0,0,0,1,0,0),
  @(0,0,0,0,0,1),
  @(0,1,0,0,0,0)
)

for ($i = 0; $i -lt 6; $i++) {
  $row = $i + 2
  for ($j = 0; $j -lt 6; $j++) {
    $col = $j + 1
    $ws.Cells.Item($row, $col).Value = $values[$i][$j]
  }
}
